$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.532.75"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.695.13"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.83"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.543"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "2.695.33"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.25"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "3.188.23"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "68.580.07"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "2.689.51"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.84"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "365.15"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.61"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.53"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.89"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.36"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "585.15"
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.95"
$ws.Range("E34").Value = "  +5.42%  "
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.86"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.78"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.50"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.93"
$ws.Range("E48").Value = "  +5.44%  "
$ws.Range("E49").Value = "  +4.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.601"
$ws.Range("E50").Value = "  +6.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.01"
$ws.Range("E51").Value = "  -0.12%  "
